$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly price record (3 rows: Pinton / Primera Maduro / Primera Pinton)
# per date, ordered chronologically by insertion. A new record for 2023-02-27 (serial 44984)
# was inserted at the top of this block (row 915), pushing the existing rows 915-1040 down
# by 3 rows (to 918-1043); this also means the previously-last record (date 44649, at what
# was rows 1038-1040) now also appears at the very end (rows 1041-1043) after the shift.
$ws.Rows("915:917").Insert()

# Row 915: new record, Fecha 2023-02-27 (serial 44984)
$ws.Cells.Item(915, 1).Value2 = 8
$ws.Cells.Item(915, 2).Value2 = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(915, 3).Value2 = 'Coquimbo'
$ws.Cells.Item(915, 4).Value2 = 44984
$ws.Cells.Item(915, 5).Value2 = 4
$ws.Cells.Item(915, 6).Value2 = 'Fruta'
$ws.Cells.Item(915, 7).Value2 = 100108
$ws.Cells.Item(915, 8).Value2 = 'Tropicales y subtropicales'
$ws.Cells.Item(915, 9).Value2 = 100108006
$ws.Cells.Item(915, 10).Value2 = 'Plátano'
$ws.Cells.Item(915, 11).Value2 = 'Sin especificar'
$ws.Cells.Item(915, 12).Value2 = 'Pintón'
$ws.Cells.Item(915, 13).Value2 = 80
$ws.Cells.Item(915, 14).Value2 = 23000
$ws.Cells.Item(915, 15).Value2 = 23000
$ws.Cells.Item(915, 16).Value2 = 23000
$ws.Cells.Item(915, 17).Value2 = '$/caja 20 kilos'
$ws.Cells.Item(915, 18).Value2 = 'Ecuador'
$ws.Cells.Item(915, 19).Value2 = 1150
$ws.Cells.Item(915, 20).Value2 = 20

# Row 916: new record, Fecha 2023-02-27 (serial 44984)
$ws.Cells.Item(916, 1).Value2 = 8
$ws.Cells.Item(916, 2).Value2 = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(916, 3).Value2 = 'Coquimbo'
$ws.Cells.Item(916, 4).Value2 = 44984
$ws.Cells.Item(916, 5).Value2 = 4
$ws.Cells.Item(916, 6).Value2 = 'Fruta'
$ws.Cells.Item(916, 7).Value2 = 100108
$ws.Cells.Item(916, 8).Value2 = 'Tropicales y subtropicales'
$ws.Cells.Item(916, 9).Value2 = 100108006
$ws.Cells.Item(916, 10).Value2 = 'Plátano'
$ws.Cells.Item(916, 11).Value2 = 'Sin especificar'
$ws.Cells.Item(916, 12).Value2 = 'Primera Maduro'
$ws.Cells.Item(916, 13).Value2 = 120
$ws.Cells.Item(916, 14).Value2 = 25000
$ws.Cells.Item(916, 15).Value2 = 25000
$ws.Cells.Item(916, 16).Value2 = 25000
$ws.Cells.Item(916, 17).Value2 = '$/caja 20 kilos'
$ws.Cells.Item(916, 18).Value2 = 'Ecuador'
$ws.Cells.Item(916, 19).Value2 = 1250
$ws.Cells.Item(916, 20).Value2 = 20

# Row 917: new record, Fecha 2023-02-27 (serial 44984)
$ws.Cells.Item(917, 1).Value2 = 8
$ws.Cells.Item(917, 2).Value2 = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(917, 3).Value2 = 'Coquimbo'
$ws.Cells.Item(917, 4).Value2 = 44984
$ws.Cells.Item(917, 5).Value2 = 4
$ws.Cells.Item(917, 6).Value2 = 'Fruta'
$ws.Cells.Item(917, 7).Value2 = 100108
$ws.Cells.Item(917, 8).Value2 = 'Tropicales y subtropicales'
$ws.Cells.Item(917, 9).Value2 = 100108006
$ws.Cells.Item(917, 10).Value2 = 'Plátano'
$ws.Cells.Item(917, 11).Value2 = 'Sin especificar'
$ws.Cells.Item(917, 12).Value2 = 'Primera Pintón'
$ws.Cells.Item(917, 13).Value2 = 120
$ws.Cells.Item(917, 14).Value2 = 26000
$ws.Cells.Item(917, 15).Value2 = 26000
$ws.Cells.Item(917, 16).Value2 = 26000
$ws.Cells.Item(917, 17).Value2 = '$/caja 20 kilos'
$ws.Cells.Item(917, 18).Value2 = 'Ecuador'
$ws.Cells.Item(917, 19).Value2 = 1300
$ws.Cells.Item(917, 20).Value2 = 20

